$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111470486
$ws.Range("B2").Value = 78578
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("Q2").Value = 554488.5866359913
$ws.Range("R2").Value = 7003175.257923778
$ws.Range("S2").Value = 22

# Row 3
$ws.Range("A3").Value = 111471083
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("Q3").Value = 554499.1143642976
$ws.Range("R3").Value = 7003141.52872613
$ws.Range("Z3").Value = "15:31"
$ws.Range("AB3").Value = "15:31"

# Row 5
$ws.Range("A5").Value = 111470792
$ws.Range("B5").Value = 96348
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("Q5").Value = 554440.9784625648
$ws.Range("R5").Value = 7003152.756292564
$ws.Range("Z5").Value = "15:19"
$ws.Range("AB5").Value = "15:19"

# Row 6
$ws.Range("A6").Value = 111470101
$ws.Range("B6").Value = 77515
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 554474.9281677724
$ws.Range("R6").Value = 7003314.266989549
$ws.Range("S6").Value = 25
$ws.Range("Z6").Value = "00:00"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AC6").Value = "På tall"

# Row 7
$ws.Range("A7").Value = 111469986
$ws.Range("Q7").Value = 554489.6113782075
$ws.Range("R7").Value = 7003329.432399829
$ws.Range("Z7").Value = "00:00"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AC7").Value = "Rikligt på tall"

# Row 8
$ws.Range("A8").Value = 111470685
$ws.Range("B8").Value = 77267
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6446
$ws.Range("F8").Value = "Kolflarnlav"
$ws.Range("G8").Value = "Carbonicola anthracophila"
$ws.Range("H8").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q8").Value = 554457.9939421143
$ws.Range("R8").Value = 7003163.892755959
$ws.Range("Z8").Value = "14:41"
$ws.Range("AB8").Value = "14:41"

# Row 9
$ws.Range("A9").Value = 111471797
$ws.Range("Q9").Value = 554597.2688619854
$ws.Range("R9").Value = 7003280.616068945
$ws.Range("Z9").Value = "15:49"
$ws.Range("AB9").Value = "15:49"
$ws.Range("AC9").Value = "På tall"

# Row 10
$ws.Range("A10").Value = 111470636
$ws.Range("B10").Value = 94134
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 53
$ws.Range("F10").Value = "Vedtrappmossa"
$ws.Range("G10").Value = "Crossocalyx hellerianus"
$ws.Range("H10").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q10").Value = 554457.9939421143
$ws.Range("R10").Value = 7003163.892755959
$ws.Range("Z10").Value = "14:41"
$ws.Range("AB10").Value = "14:41"

# Row 11
$ws.Range("A11").Value = 111470743
$ws.Range("B11").Value = 78611
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 6463
$ws.Range("F11").Value = "Bårdlav"
$ws.Range("G11").Value = "Nephroma parile"
$ws.Range("Q11").Value = 554457.9939421143
$ws.Range("R11").Value = 7003163.892755959
$ws.Range("Z11").Value = "14:41"
$ws.Range("AB11").Value = "14:41"
$ws.Range("AC11").ClearContents()

# Row 13
$ws.Range("A13").Value = 111471685
$ws.Range("B13").Value = 77515
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("Q13").Value = 554595.0694405095
$ws.Range("R13").Value = 7003142.694495555
$ws.Range("Z13").Value = "15:49"
$ws.Range("AB13").Value = "15:49"
